$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add C1, D1, E1 (B1 already has 16)
$ws.Range("C1").Value = 17
$ws.Range("D1").Value = 18
$ws.Range("E1").Value = 19

# Column A: rows 4-9
$ws.Range("A4").Value = 8
$ws.Range("A5").Value = 9
$ws.Range("A6").Value = 10
$ws.Range("A7").Value = 11
$ws.Range("A8").Value = 12
$ws.Range("A9").Value = 13
